$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.899.59'
$ws.Range('E2').Value = '  +2.40%  '
$ws.Range('D3').Value = '2.455.71'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.42'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.19'
$ws.Range('E6').Value = '  +9.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.507'
$ws.Range('E8').Value = '  +2.29%  '
$ws.Range('D9').Value = '2.465.01'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.78'
$ws.Range('E10').Value = '  +8.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.330'
$ws.Range('E12').Value = '  +4.04%  '
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').Value = '2.883.22'
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').Value = '56.097.43'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.98'
$ws.Range('E16').Value = '  +5.80%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '2.470.33'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.49'
$ws.Range('E19').Value = '  +6.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.02'
$ws.Range('E20').Value = '  +4.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '315.95'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  +0.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.75'
$ws.Range('E23').Value = '  +6.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '58.25'
$ws.Range('E24').Value = '  +2.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.411'
$ws.Range('E25').Value = '  +5.98%  '
$ws.Range('E26').Value = '  -0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.159'
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').Value = '2.594.65'
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.63'
$ws.Range('E29').Value = '  +6.15%  '
$ws.Range('D30').Value = '0.0₃0777'
$ws.Range('E30').Value = '  +7.32%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '147.92'
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.14'
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('E34').Value = '  +3.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.16'
$ws.Range('E35').Value = '  +2.23%  '
$ws.Range('E36').Value = '  +7.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.71'
$ws.Range('E37').Value = '  +3.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.859'
$ws.Range('E38').Value = '  +6.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.88'
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.50'
$ws.Range('E40').Value = '  +7.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0552'
$ws.Range('E42').Value = '  +4.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.601'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  +5.83%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '259.61'
$ws.Range('E45').Value = '  +10.17%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.70'
$ws.Range('E46').Value = '  +10.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0917'
$ws.Range('E47').Value = '  +3.88%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.19'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0227'
$ws.Range('E49').Value = '  +3.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.45'
$ws.Range('E50').Value = '  +4.73%  '
$ws.Range('D51').Value = '1.862.40'
$ws.Range('E51').Value = '  -4.12%  '
